# Generate Report for Handoff
# Refreshes the localization-status report: new handoff GUID, new target-file
# hash, and updated handoff/handback timestamps across the Overview, zh-cn
# and de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "45c3539c-16a1-40f7-a8e0-c3bc605a1c64"
$newGuid = "ae3582fa-a85f-4867-85d9-f0eb8e430ccc"
$oldHash = "fb0c7af5c9185ae917426c4ef0316033e593ea50"
$newHash = "a84f5a4f2f3be76bee201d47c3bb2686c6ee6f93"

$oldMd = "$oldGuid.md"
$newMd = "$newGuid.md"

$oldZh = "$oldGuid.$oldHash.zh-cn.xlf"
$newZh = "$newGuid.$newHash.zh-cn.xlf"

$oldDe = "$oldGuid.$oldHash.de-de.xlf"
$newDe = "$newGuid.$newHash.de-de.xlf"

# ---- Overview sheet ----
$ws1 = $wb.Worksheets.Item("Overview")
$link1 = "https://github.com/OpenLocalizationTest/oltest/blob/62d070b725ad63393d4e2f0bfa9db21a6ac8cde4/e2e/$oldMd"

$ws1.Hyperlinks.Delete()
$ws1.Range("A2").Value = $newMd
$ws1.Range("D2").Value = "2016-49-20 20:49:22"
$ws1.Hyperlinks.Add($ws1.Range("A2"), $link1, "", "", $newMd)

# ---- zh-cn sheet ----
$ws2 = $wb.Worksheets.Item("zh-cn")
$linkMd2 = "https://github.com/OpenLocalizationTest/oltest/blob/62d070b725ad63393d4e2f0bfa9db21a6ac8cde4/e2e/$oldMd"
$linkZh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/76cdc180b560c3de374063f69d19d4f931093dbd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldZh"

$ws2.Hyperlinks.Delete()
$ws2.Range("A2").Value = $newMd
$ws2.Range("D2").Value = $newZh
$ws2.Range("E2").Value = "2016-03-20 20:49:19"
$ws2.Hyperlinks.Add($ws2.Range("A2"), $linkMd2, "", "", $newMd)
$ws2.Hyperlinks.Add($ws2.Range("B2"), $linkMd2, "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), $linkZh, "", "", $newZh)

# ---- de-de sheet ----
$ws3 = $wb.Worksheets.Item("de-de")
$linkMd3 = "https://github.com/OpenLocalizationTest/oltest/blob/62d070b725ad63393d4e2f0bfa9db21a6ac8cde4/e2e/$oldMd"
$linkDe = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2f381aa8ebe340e5397e02c5683d2e395cbf689c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldDe"

$ws3.Hyperlinks.Delete()
$ws3.Range("A2").Value = $newMd
$ws3.Range("D2").Value = $newDe
$ws3.Range("E2").Value = "2016-03-20 20:49:22"
$ws3.Hyperlinks.Add($ws3.Range("A2"), $linkMd3, "", "", $newMd)
$ws3.Hyperlinks.Add($ws3.Range("B2"), $linkMd3, "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), $linkDe, "", "", $newDe)
